# Update recurrence metrics for the most recent quarters (row 19 = 2025Q1, row 20 = 2025Q2)
# per commit "atualizei dados para BIBI e ADD 06-05-2025"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (2025Q1)
$ws.Range("C19").Value = 324
$ws.Range("D19").Value = 263
$ws.Range("E19").Value = 61
$ws.Range("F19").Value = 82.44514106583071

# Row 20 (2025Q2)
$ws.Range("C20").Value = 235
$ws.Range("D20").Value = 207
$ws.Range("E20").Value = 28
$ws.Range("F20").Value = 63.88888888888889
